# Update "Förändrad" (C column) dates from 45205 to 45206 for existing rows (2..518)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C518").Value = 45206

# Force row 518 to carry an explicit custom height (matches target diff: ht="15" customHeight="1")
$ws.Rows.Item(518).RowHeight = 15

# Add new row 519 with the new avverkningsanmälan entry
$ws.Range("A519").Value = "A 48120-2023"

$ws.Range("B519").Value = 45205
$ws.Range("B519").NumberFormat = "YYYY-MM-DD"

$ws.Range("C519").Value = 45206
$ws.Range("C519").NumberFormat = "YYYY-MM-DD"

$ws.Range("D519").Value = "HALLANDS LÄN"
$ws.Range("E519").Value = "FALKENBERG"

$ws.Range("G519").Value = 3
$ws.Range("H519:Q519").Value = 0

$ws.Range("R519").Value = ""
$ws.Range("R519").WrapText = $true
